# Insert two new price-report rows (100 and 101) above the current row 100
# in the "Pepino ensalada" consolidated sheet. Everything that was at row
# 100 onward shifts down by two rows; the new rows carry a fresh weekly
# observation (fecha 44588) for "Primera" and "Segunda" calidad.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 100 (rows 100:101), pushing old data down.
$ws.Range("A100:A101").EntireRow.Insert()

# New row 100 - Primera
$ws.Cells.Item(100, 1).Value  = 2
$ws.Cells.Item(100, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(100, 3).Value  = "Coquimbo"
$ws.Cells.Item(100, 4).Value  = 44588
$ws.Cells.Item(100, 5).Value  = 4
$ws.Cells.Item(100, 6).Value  = 100112043
$ws.Cells.Item(100, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(100, 8).Value  = "Sin especificar"
$ws.Cells.Item(100, 9).Value  = "Primera"
$ws.Cells.Item(100, 10).Value = 480
$ws.Cells.Item(100, 11).Value = 9000
$ws.Cells.Item(100, 12).Value = 10000
$ws.Cells.Item(100, 13).Value = 9417
$ws.Cells.Item(100, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(100, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(100, 16).Value = 135
$ws.Cells.Item(100, 17).Value = 70
$ws.Cells.Item(100, 18).Value = "Hortaliza"

# New row 101 - Segunda
$ws.Cells.Item(101, 1).Value  = 2
$ws.Cells.Item(101, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(101, 3).Value  = "Coquimbo"
$ws.Cells.Item(101, 4).Value  = 44588
$ws.Cells.Item(101, 5).Value  = 4
$ws.Cells.Item(101, 6).Value  = 100112043
$ws.Cells.Item(101, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(101, 8).Value  = "Sin especificar"
$ws.Cells.Item(101, 9).Value  = "Segunda"
$ws.Cells.Item(101, 10).Value = 300
$ws.Cells.Item(101, 11).Value = 7000
$ws.Cells.Item(101, 12).Value = 8000
$ws.Cells.Item(101, 13).Value = 7500
$ws.Cells.Item(101, 14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(101, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(101, 16).Value = 75
$ws.Cells.Item(101, 17).Value = 100
$ws.Cells.Item(101, 18).Value = "Hortaliza"
